$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card21")

# Helper: force a literal text value (leading apostrophe = Excel's
# "treat as text" quote-prefix) then strip the quote-prefix styling so the
# cell ends up as a plain, unstyled inline/shared string - matching how the
# rest of this sheet's text cells are stored.
function Set-TextValue {
    param($range, $value)
    $range.Value = "'" + $value
    $range.Style = "Normal"
}

# Row 30: fill previously-empty cells with literal "nan" text (matches the
# sheet's convention for missing values elsewhere in the table)
foreach ($col in @("D30","E30","F30","G30","I30","J30","K30","M30","N30")) {
    Set-TextValue $ws.Range($col) "nan"
}

# Row 31: new service-log entry for Card21
Set-TextValue $ws.Range("A31") "21"
foreach ($col in @("B31","C31","D31","E31","F31","G31","H31","I31","J31","K31")) {
    Set-TextValue $ws.Range($col) ""
}
Set-TextValue $ws.Range("L31") "21/12/2025"
Set-TextValue $ws.Range("M31") "خلل ف جوده وحدوثt-cone"
Set-TextValue $ws.Range("N31") "تم اعاده عيار الماكينه وتم تغير جرائد اماميه 550"
Set-TextValue $ws.Range("O31") "م.شحته ،تيم الكرد"
